# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals.
# Column G ("K") holds recomputed values; update each affected row's K value
# to its freshly-calculated figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2=1; 3=1; 4=1; 5=1; 6=2; 7=0; 8=1; 9=1; 10=2;
    12=1; 13=1; 14=1; 15=1; 16=1; 17=0; 18=0; 19=2; 20=0; 21=3; 22=0; 23=0; 24=1; 25=0; 26=0; 27=1; 28=1; 29=1; 30=4; 31=1; 32=0; 33=2; 34=1; 35=0; 36=1; 37=0; 38=2; 39=2; 40=1;
    42=0; 43=1; 44=0; 45=2; 46=1; 47=1; 48=1; 49=1; 50=1; 51=1; 52=1; 53=1; 54=3;
    56=0; 57=2; 58=0; 59=2; 60=1; 61=1; 62=0; 63=0; 64=0; 65=1; 66=3; 67=0; 68=0; 69=1;
    71=2; 72=1; 73=2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Output "Updated $($kValues.Keys.Count) K column values"
